$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (targetX) values from 7/-7 to 300/-300 per the diff
$ws.Range("C2").Value = 300
$ws.Range("C3").Value = 300
$ws.Range("C4").Value = 300
$ws.Range("C5").Value = 300
$ws.Range("C6").Value = -300
$ws.Range("C7").Value = -300
$ws.Range("C8").Value = -300
$ws.Range("C9").Value = -300
$ws.Range("C10").Value = -300
$ws.Range("C11").Value = -300
$ws.Range("C12").Value = -300
$ws.Range("C13").Value = 300
$ws.Range("C14").Value = -300
$ws.Range("C15").Value = 300

# Update column B (cueOri) values from -180 to 180 for rows 12-15
$ws.Range("B12").Value = 180
$ws.Range("B13").Value = 180
$ws.Range("B14").Value = 180
$ws.Range("B15").Value = 180

# Update the active cell selection from H20 to H21
$ws.Range("H21").Select()
